$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Highlight the "done" Wednesday-column cells (and the Review Component
#     Access cell) with the existing yellow fill used elsewhere in the sheet ---
$ws.Range("E2").Interior.Color = 65535
$ws.Range("C3").Interior.Color = 65535
$ws.Range("E4").Interior.Color = 65535
$ws.Range("E5").Interior.Color = 65535

# --- Add the new Zombie Spawner follow-up task as row 8 ---
$ws.Range("E8").Value = "If done, make plan for Thursday"
$ws.Range("E8").HorizontalAlignment = -4108
$ws.Range("E8").WrapText = $true
$ws.Rows.Item(8).RowHeight = 30

# --- Move the active selection to E5 ---
$ws.Range("E5").Select() | Out-Null
